$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 39
$ws.Range("I2").Value = 109
$ws.Range("J2").Value = 523
$ws.Range("K2").Value = 8
$ws.Range("L2").Value = 140
$ws.Range("M2").Value = 9
$ws.Range("N2").Value = 98
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 12
$ws.Range("S2").Value = 53
$ws.Range("T2").Value = 87
$ws.Range("U2").Value = 7
$ws.Range("V2").Value = 878
$ws.Range("W2").Value = 1
$ws.Range("X2").Value = 866
$ws.Range("Y2").Value = 3
